# Generate Report for Handoff
#
# This mirrors the behavior of the localization-status report generator:
#  - For files whose handoff to the zh-cn / de-de xliff targets just completed,
#    the "Priority" column (E) on each language sheet is stamped with the
#    handoff type "ht".
#  - The "Latest Handoff Datetime" (H) on each language sheet, and the
#    "Latest HO Xliff Generate Date" (G) on the Overview sheet, are refreshed
#    to the new generation timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows (on the zh-cn / de-de tables, and matching rows on Overview) whose
# handoff just completed for this generation pass.
$rows = @(8, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Mark the handoff type as "ht" (Priority column) on both language sheets.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # Refresh the "Latest Handoff Datetime" for each language sheet.
    $wsZhCn.Range("H$r").Value = "2016-08-26 12:20:58"
    $wsDeDe.Range("H$r").Value = "2016-08-26 12:21:09"

    # Refresh the overall "Latest HO Xliff Generate Date" on the Overview sheet
    # (the most recent of the per-language handoff times).
    $wsOverview.Range("G$r").Value = "2016-08-26 12:21:09"
}
